$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Data for new rows 191-211 (SXT antibiotic model), appended below existing row 190 data.
# Each entry: row number, then values for columns A..H (or $null for blank/NA cells).
$newRows = @(
    ,("(Intercept)", [double]"1.949993489118026e-41", [double]"19.26715602857632", [double]"-4.865178993727496", [double]"1.143532694916006e-06", [double]"7.524190342496107e-58", [double]"4.786175685614378e-25", "SXT")
    ,("Year", [double]"1.047992406395753", [double]"0.009550447799722394", [double]"4.908287134806678", [double]"9.187527592856661e-07", [double]"1.028570266527557", [double]"1.067810524035388", "SXT")
    ,("Specimen_typeother", [double]"0.9990288573289787", [double]"0.08591900568955699", [double]"-0.01130849371207797", [double]"0.9909773197680506", [double]"0.8441719451610779", [double]"1.182280513770871", "SXT")
    ,("Specimen_typeRespiratory", [double]"0.653224920749101", [double]"0.09558079878319813", [double]"-4.455222934364498", [double]"8.380606335151297e-06", [double]"0.5414153342584654", [double]"0.7875515843106367", "SXT")
    ,("Specimen_typeUrine", [double]"0.871934037114984", [double]"0.07404593993396408", [double]"-1.850763236539096", [double]"0.0642036221377104", [double]"0.7542871774363854", [double]"1.008350330457849", "SXT")
    ,("Specimen_typeWound & soft tissues", [double]"1.000668250665032", [double]"0.08401908789825828", [double]"0.007950901416432178", [double]"0.9936561653549565", [double]"0.8487283628100275", [double]"1.179838262617296", "SXT")
    ,("HospitalCHBH", [double]"0.3875487016949927", [double]"0.1600587204351798", [double]"-5.922287481214648", [double]"3.174939657609115e-09", [double]"0.2826074560157678", [double]"0.5294255678109484", "SXT")
    ,("HospitalCNGMO", [double]"1.30107387830818", [double]"0.4753114263927906", [double]"0.5537211375328404", [double]"0.5797697069592305", [double]"0.5290678256955028", [double]"3.510129777184626", "SXT")
    ,("HospitalRabta", [double]"0.2989064061391349", [double]"0.1937184116410836", [double]"-6.233918435057435", [double]"4.549100063747805e-10", [double]"0.2037614931454496", [double]"0.4357156652452083", "SXT")
    ,("HospitalTCB", [double]"0.9313338815755413", [double]"0.1344625481022051", [double]"-0.5290502089522586", [double]"0.5967706202619545", [double]"0.7144137173457489", [double]"1.210586364456798", "SXT")
    ,("Ward_ED_ICUED", [double]"0.2042614581068638", [double]"0.165887934726535", [double]"-9.574864205714878", [double]"1.019912466829873e-21", [double]"0.1471148733333632", [double]"0.2819680781964216", "SXT")
    ,("Ward_ED_ICUOther", [double]"0.2781029051115721", [double]"0.1245347088180019", [double]"-10.276364586704", [double]"9.006079417170985e-25", [double]"0.2173875787481296", [double]"0.3543194910725941", "SXT")
    ,("GenderF", [double]"0.7364947102219679", [double]"0.04667786420796057", [double]"-6.552425450781537", [double]"5.66099816982559e-11", [double]"0.6720719930984718", [double]"0.8070214661807554", "SXT")
    ,("HospitalCHBH:Ward_ED_ICUED", [double]"2.343788196993692", [double]"0.2857815781916495", [double]"2.980487802529208", [double]"0.002877896989144335", [double]"1.329218476924106", [double]"4.08188467832268", "SXT")
    ,("HospitalCNGMO:Ward_ED_ICUED", $null, $null, $null, $null, $null, $null, "SXT")
    ,("HospitalRabta:Ward_ED_ICUED", [double]"2.338670033406201", [double]"0.2798031821250126", [double]"3.036357197749866", [double]"0.002394555052029372", [double]"1.348618240176538", [double]"4.042353319373585", "SXT")
    ,("HospitalTCB:Ward_ED_ICUED", [double]"1.303485683578708", [double]"0.4313454542708406", [double]"0.6144540730518631", [double]"0.5389153131547566", [double]"0.541652139900454", [double]"2.983176629121618", "SXT")
    ,("HospitalCHBH:Ward_ED_ICUOther", [double]"3.001174941099489", [double]"0.1761185766910956", [double]"6.240135933844147", [double]"4.371907995547408e-10", [double]"2.128088385240857", [double]"4.245445594512282", "SXT")
    ,("HospitalCNGMO:Ward_ED_ICUOther", [double]"2.012433320594806", [double]"0.5109833199503171", [double]"1.368625099555075", [double]"0.1711164940604356", [double]"0.7015380192432163", [double]"5.329289847359406", "SXT")
    ,("HospitalRabta:Ward_ED_ICUOther", [double]"2.895238107311508", [double]"0.211369827567489", [double]"5.029418662317643", [double]"4.919691389254079e-07", [double]"1.917344812397091", [double]"4.393179137146878", "SXT")
    ,("HospitalTCB:Ward_ED_ICUOther", [double]"1.99619346290885", [double]"0.1547658455126479", [double]"4.466373677006553", [double]"7.955660186636874e-06", [double]"1.475382356051192", [double]"2.706917107712356", "SXT")
)

$startRow = 191
for ($i = 0; $i -lt $newRows.Count; $i++) {
    $r = $startRow + $i
    $rowVals = $newRows[$i]
    for ($c = 1; $c -le 8; $c++) {
        $val = $rowVals[$c - 1]
        if ($null -ne $val) {
            $ws.Cells.Item($r, $c).Value = $val
        }
    }
}
